$wb = $excel.ActiveWorkbook

$sheetNames = @("Sheet1", "DeviceSetupLogins", "InjectSpecificUser")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Update passwords for Hasanein.Ba-Alawi.* rows (8,9,10): Appian1 -> Appian2
    $ws.Range("B8").Value = "Appian2"
    $ws.Range("B9").Value = "Appian2"
    $ws.Range("B10").Value = "Appian2"

    # Update passwords for Yaaseen.Choudhury.Business/.AuthorisedRep (17,19): password2 -> MHRA12345
    $ws.Range("B17").Value = "MHRA12345"
    $ws.Range("B19").Value = "MHRA12345"

    # Update the selection to B17:B19 with active cell B17
    $ws.Range("B17:B19").Select()
}

# Re-select the InjectSpecificUser sheet (tabSelected) and restore its selection
$ws3 = $wb.Worksheets.Item("InjectSpecificUser")
$ws3.Activate()
$ws3.Range("B17:B19").Select()
